$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '40.800.72'
$ws.Range('E2').Value = '  +3.61%  '

# Row 3
$ws.Range('D3').Value = '2.214.28'
$ws.Range('E3').Value = '  +2.46%  '

# Row 4
$ws.Range('E4').Value = '  +0.11%  '

# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '229.11'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.42%  '

# Row 6
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.632'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +1.68%  '

# Row 7
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '64.51'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +0.38%  '

# Row 8
$ws.Range('E8').Value = '  +0.09%  '

# Row 9
$ws.Range('E9').Value = '  +1.71%  '

# Row 10
$ws.Range('E10').Value = '  +0.86%  '

# Row 11
$ws.Range('E11').Value = '  -0.07%  '

# Row 12
$ws.Range('D12').Value = '2.541.89'
$ws.Range('E12').Value = '  +2.50%  '

# Row 13
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '15.90'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -0.47%  '

# Row 14
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '22.18'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -0.51%  '

# Row 15
$ws.Range('E15').Value = '  +0.80%  '

# Row 16
$ws.Range('E16').Value = '  +1.12%  '

# Row 17
$ws.Range('D17').Value = '2.218.65'
$ws.Range('E17').Value = '  +2.99%  '

# Row 18
$ws.Range('D18').Value = '40.735.27'
$ws.Range('E18').Value = '  +3.47%  '

# Row 19
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '73.89'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +2.78%  '

# Row 20
$ws.Range('B20').Value = 'ShibaInu'
$ws.Range('C20').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D20').Value = '0.0₃0902'
$ws.Range('E20').Value = '  +5.74%  '

# Row 21
$ws.Range('B21').Value = 'Uniswap'
$ws.Range('C21').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '6.15'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +0.20%  '

# Row 22
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '250.00'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +8.14%  '

# Row 23
$ws.Range('E23').Value = '  +0.00%  '

# Row 24
$ws.Range('E24').Value = '  +1.22%  '

# Row 25
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.30'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -9.67%  '

# Row 26
$ws.Range('B26').Value = 'Cosmos'
$ws.Range('C26').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '9.71'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +0.56%  '

# Row 27
$ws.Range('B27').Value = 'Monero'
$ws.Range('C27').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '172.92'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +0.27%  '

# Row 28
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.145'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +4.13%  '

# Row 29
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '20.32'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +2.04%  '

# Row 30
$ws.Range('E30').Value = '  +2.08%  '

# Row 31
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '2.82'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +8.02%  '

# Row 32
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.123'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +1.23%  '

# Row 33
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '4.66'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +0.82%  '

# Row 34
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '7.15'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +0.69%  '

# Row 35
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '4.77'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -0.46%  '

# Row 36
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.0631'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +2.05%  '

# Row 37
$ws.Range('E37').Value = '  +6.47%  '

# Row 38
$ws.Range('E38').Value = '  +1.81%  '

# Row 39
$ws.Range('E39').Value = '  +0.07%  '

# Row 40
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '4.93'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +14.39%  '

# Row 41
$ws.Range('E41').Value = '  +1.32%  '

# Row 42
$ws.Range('E42').Value = '  +8.67%  '

# Row 43
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '101.41'
$ws.Range('D43').Style = 'Normal'

# Row 44
$ws.Range('E44').Value = '  +4.47%  '

# Row 45
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '17.41'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -2.52%  '

# Row 46
$ws.Range('D46').Value = '1.510.26'
$ws.Range('E46').Value = '  -1.86%  '

# Row 47
$ws.Range('B47').Value = 'HuobiToken'
$ws.Range('C47').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.89'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +2.49%  '

# Row 48
$ws.Range('B48').Value = 'Cronos'
$ws.Range('C48').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.0939'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +1.58%  '

# Row 49
$ws.Range('B49').Value = 'ARBITRUM'
$ws.Range('C49').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.11'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +0.81%  '

# Row 50
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.000207'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +37.71%  '

# Row 51
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '9.56'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +10.08%  '
